# Fix image proportions on slide 12 ("Casi interessanti - 3")
# Restores the two pictures to their native aspect ratio while keeping them
# roughly centered in the same place (author commit: "fixed images
# proportion slide casi particolari 3").
#
# Left/Top/Width/Height are expressed in points; 1 pt = 12700 EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# Picture "Immagine 4" (left picture)
$pic1 = $s.Shapes.Item("Immagine 4")
$pic1.Left   = 125.46748031496062
$pic1.Top    = 224.64063292125982
$pic1.Width  = 263.6888278976378
$pic1.Height = 270.39283764566926

# Picture "Immagine 6" (right picture)
$pic2 = $s.Shapes.Item("Immagine 6")
$pic2.Left   = 411.8359842519685
$pic2.Top    = 265.67299212598425
$pic2.Width  = 510.2929233858267
$pic2.Height = 185.76448918897637
